$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the Google Colab spell-check run-splits (proofErr wrappers)
#    by re-finding/replacing the sentences verbatim -- this merges the
#    split runs back into one run and drops the now-pointless proofErr
#    start/end markers.
# ---------------------------------------------------------------------

$t1 = "The artifact is based on Google Colab computing engine. To run a certain version of the artifact, you"
$d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2) | Out-Null

$t2 = "Sign-in to Google Colab using a google account"
$d.Content.Find.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2) | Out-Null

$t3 = "Upload all the dependent module to Google Colab “Files”"
$d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2) | Out-Null

# ---------------------------------------------------------------------
# 2. "Open one of the artifact versions on Google Colab" paragraph --
#    here the proofErr "spellEnd" marker sits as the very last child of
#    the paragraph (nothing after it), so a plain Find/Replace leaves it
#    orphaned. Work around it by appending a throw-away marker character
#    first (so the proofErr is no longer paragraph-final), then running
#    the real replace across that marker -- this consolidates the run
#    and drops the proofErr cleanly.
# ---------------------------------------------------------------------

$p5 = $d.Paragraphs.Item(6)
$p5end = $p5.Range.End - 1
$marker = $d.Range($p5end, $p5end)
$marker.InsertAfter([char]1)

$t4 = "Open one of the artifact versions on Google Colab" + [char]1
$t4new = "Open one of the artifact versions on Google Colab"
$d.Content.Find.Execute($t4, $true, $false, $false, $false, $false, $true, 1, $false, $t4new, 2) | Out-Null

# ---------------------------------------------------------------------
# 3. "Colab Notebooks" inside the "Mount your google drive..." paragraph
# ---------------------------------------------------------------------

$t5 = "“Colab Notebooks”"
$d.Content.Find.Execute($t5, $true, $false, $false, $false, $false, $true, 1, $false, $t5, 2) | Out-Null

# ---------------------------------------------------------------------
# 4. Rewrite the tail of the "Mount your google drive..." paragraph to
#    describe the new per-dataset subfolder layout.
# ---------------------------------------------------------------------

$oldTail = "”. “SIT723” folder should be place where you save/store all the datasets"
$newTail = "”. Inside “SIT723” folder, you create 3 subfolders for each dataset, namely “FER2013”, “CKplus”, and “JAFFE”. Each subfolder should be the location where you save/store all the datasets according their categories."
$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

# ---------------------------------------------------------------------
# 5. "For example:" -> "For example, the training images and labels of
#    FER2013 should look like:"
# ---------------------------------------------------------------------

$d.Content.Find.Execute("For example:", $true, $false, $false, $false, $false, $true, 1, $false, "For example, the training images and labels of FER2013 should look like:", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Update the x_train.csv / y_train.csv example paths: drop the
#    proofErr wrapper around "Colab", insert the new "/FER2013" path
#    segment, and wrap the whole path in curly quotes. The proofErr
#    "spellStart" marker is the very first child of each of these two
#    paragraphs (nothing before it), so -- same issue as step 2, mirrored
#    -- a plain Find/Replace leaves it orphaned. Prime a throw-away
#    marker character immediately before the paragraph and fold it into
#    the search/replace so the whole paragraph collapses into one clean
#    run.
# ---------------------------------------------------------------------

$pX = $d.Paragraphs.Item(9)
$pXstart = $pX.Range.Start
$markerX = $d.Range($pXstart, $pXstart)
$markerX.InsertBefore([char]1)

$oldX = [char]1 + "Colab Notebooks/SIT723/x_train.csv"
$newX = [char]1 + "“Colab Notebooks/SIT723/FER2013/x_train.csv”"
$d.Content.Find.Execute($oldX, $true, $false, $false, $false, $false, $true, 1, $false, $newX, 2) | Out-Null

$pY = $d.Paragraphs.Item(10)
$pYstart = $pY.Range.Start
$markerY = $d.Range($pYstart, $pYstart)
$markerY.InsertBefore([char]1)

$oldY = [char]1 + "Colab Notebooks/SIT723/y_train.csv"
$newY = [char]1 + "“Colab Notebooks/SIT723/FER2013/y_train.csv”"
$d.Content.Find.Execute($oldY, $true, $false, $false, $false, $false, $true, 1, $false, $newY, 2) | Out-Null

# ---------------------------------------------------------------------
# 7. Remove the now-unneeded "... " placeholder paragraph entirely
#    (Expand to the whole paragraph -- including its paragraph mark --
#    before deleting, so the following paragraph collapses up in its
#    place instead of leaving a blank paragraph behind).
# ---------------------------------------------------------------------

$rngDots = $d.Content
$foundDots = $rngDots.Find.Execute("… ")
if ($foundDots) {
    $rngDots.Expand(4) | Out-Null
    $rngDots.Delete()
}

Write-Output "done"
